# Added more APIs from GT team development
# Adds a new "ArchiveSample" worksheet (cloned from WorklistOverview, the
# existing sheet with the same EndPoint/gridName/searchableColumns shape)
# and points its searchableColumns values at the new Archive Sample API.

$wb = $excel.ActiveWorkbook
$source = $wb.Worksheets.Item("WorklistOverview")

# Duplicate the WorklistOverview sheet (same row/col layout, styles and
# merged cells) and drop it right after the last existing tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$source.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "ArchiveSample"

# Re-create the merged header bands in row order so they come out sorted.
$newSheet.Range("A1:B1").UnMerge() | Out-Null
$newSheet.Range("A6:B6").UnMerge() | Out-Null
$newSheet.Range("A10:B10").UnMerge() | Out-Null
$newSheet.Range("A1:B1").Merge() | Out-Null
$newSheet.Range("A6:B6").Merge() | Out-Null
$newSheet.Range("A10:B10").Merge() | Out-Null

# Point the searchableColumns gridName values at the new Archive Sample API.
$newSheet.Range("B3").Value = "ArchiveSamples"
$newSheet.Range("B8").Value = "ArchiveSamples"
$newSheet.Range("B12").Value = "Abc"

# The previously-active WorklistOverview sheet keeps A1:B12 selected instead.
$source.Activate() | Out-Null
$source.Range("A1:B12").Select() | Out-Null

# Leave the new sheet as the active tab with G6 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("G6").Select() | Out-Null

Write-Output "Added ArchiveSample sheet"
